$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 327, shifting existing rows 327-352 down to 328-353.
$ws.Rows.Item(327).Insert()

# Populate the newly inserted row 327 with its data.
$ws.Cells.Item(327, 1).Value = 6
$ws.Cells.Item(327, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(327, 3).Value = "Metropolitana"
$ws.Cells.Item(327, 4).Value = 45166
$ws.Cells.Item(327, 5).Value = 13
$ws.Cells.Item(327, 6).Value = 100112001
$ws.Cells.Item(327, 7).Value = "Berenjena"
$ws.Cells.Item(327, 8).Value = "Sin especificar"
$ws.Cells.Item(327, 9).Value = "Primera"
$ws.Cells.Item(327, 10).Value = 220
$ws.Cells.Item(327, 11).Value = 9000
$ws.Cells.Item(327, 12).Value = 10000
$ws.Cells.Item(327, 13).Value = 9455
$ws.Cells.Item(327, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(327, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(327, 16).Value = 189
$ws.Cells.Item(327, 17).Value = 50
$ws.Cells.Item(327, 18).Value = "Hortaliza"
